$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

# Rename the "status" column (header + table column) to "state"
$ws.Range("C1").Value = "state"

# Update the mgmt_ip values for the existing switches
$ws.Range("B2").Value = "192.168.1.101"
$ws.Range("B3").Value = "192.168.1.102"

# Keep the table's column name in sync with the header cell
$tbl = $ws.ListObjects.Item("devices")
$tbl.ListColumns.Item("status").Name = "state"

# Move the active selection to B4
$ws.Range("B4").Select()
